# Weekly update: insert 3 new price records after row 204, pushing the
# existing rows 205:300 down to 208:303, then populate the 3 new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows starting at row 205 (shifts rows 205:300 -> 208:303)
$ws.Range("A205:R207").Insert()

# Copy the date number format (style) used by column D down into the newly
# inserted D205:D207 cells so they keep the same date formatting.
$ws.Range("D208").Copy()
$ws.Range("D205:D207").PasteSpecial(-4122)  # xlPasteFormats

# Row 205
$ws.Range("A205").Value = 4
$ws.Range("B205").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C205").Value = "Los Lagos"
$ws.Range("D205").Value = 44523
$ws.Range("E205").Value = 10
$ws.Range("F205").Value = 100112006
$ws.Range("G205").Value = "Repollo"
$ws.Range("H205").Value = "Copenhague"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 500
$ws.Range("K205").Value = 1300
$ws.Range("L205").Value = 1400
$ws.Range("M205").Value = 1350
$ws.Range("N205").Value = "`$/unidad"
$ws.Range("O205").Value = "Región Metropolitana"
$ws.Range("P205").Value = 1350
$ws.Range("Q205").Value = 1
$ws.Range("R205").Value = "Hortaliza"

# Row 206
$ws.Range("A206").Value = 4
$ws.Range("B206").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C206").Value = "Los Lagos"
$ws.Range("D206").Value = 44523
$ws.Range("E206").Value = 10
$ws.Range("F206").Value = 100112006
$ws.Range("G206").Value = "Repollo"
$ws.Range("H206").Value = "Crespo record"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 1000
$ws.Range("K206").Value = 1100
$ws.Range("L206").Value = 1200
$ws.Range("M206").Value = 1150
$ws.Range("N206").Value = "`$/unidad"
$ws.Range("O206").Value = "Región Metropolitana"
$ws.Range("P206").Value = 1150
$ws.Range("Q206").Value = 1
$ws.Range("R206").Value = "Hortaliza"

# Row 207
$ws.Range("A207").Value = 4
$ws.Range("B207").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C207").Value = "Los Lagos"
$ws.Range("D207").Value = 44523
$ws.Range("E207").Value = 10
$ws.Range("F207").Value = 100112006
$ws.Range("G207").Value = "Repollo"
$ws.Range("H207").Value = "Crespo record"
$ws.Range("I207").Value = "Segunda"
$ws.Range("J207").Value = 500
$ws.Range("K207").Value = 1000
$ws.Range("L207").Value = 1000
$ws.Range("M207").Value = 1000
$ws.Range("N207").Value = "`$/unidad"
$ws.Range("O207").Value = "Región Metropolitana"
$ws.Range("P207").Value = 1000
$ws.Range("Q207").Value = 1
$ws.Range("R207").Value = "Hortaliza"
